$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / selection state -------------------------------------------------
# Best-effort: persist the new window height if the host supports it (some
# hosts don't round-trip workbookView geometry through the COM window object).
try {
    $win = $wb.Windows.Item(1)
    $win.Height = 16440
} catch {}

# Selection moves from F12 to the merged header cell J1 (sqref J1:L1)
[void]$ws.Range("J1:L1").Select()

# --- Data edits ----------------------------------------------------------------
# Row 3 (HGL Temperature Rise) - RP all (new McT algorithm) columns J:L
$ws.Range("J3").Value = 1.1100000000000001
$ws.Range("K3").Formula = "=0.44/2"

# Row 4 (HGL Depth)
$ws.Range("J4").Value = 1.01
$ws.Range("K4").Formula = "=0.32/2"

# Row 5 (Ceiling Jet Temp. Rise)
$ws.Range("J5").Value = 1.25
$ws.Range("K5").Formula = "=0.53/2"

# Row 6 (Plume Temperature Rise)
$ws.Range("K6").Formula = "=0.42/2"

# Row 8 (Oxygen Concentration)
$ws.Range("J8").Value = 1.03
$ws.Range("K8").Formula = "=0.63/2"

# Row 9 (Carbon Dioxide Concentration)
$ws.Range("K9").Formula = "=0.56/2"

# Row 10 (Smoke Concentration) - K10 becomes blank (was a literal 0)
$ws.Range("K10").ClearContents()

# Row 11 (Room Pressure Rise) - K11 becomes blank (was a literal 0)
$ws.Range("K11").ClearContents()

# Row 13 (Radiant Heat Flux)
$ws.Range("K13").Formula = "=1.29/2"

# Row 15 (Wall Temperature Rise)
$ws.Range("J15").Value = 0.99
$ws.Range("K15").Formula = "=0.99/2"
